# chore: adapt column header formatting to respective input file names (#7)
#
# The AHB-diff sheet had its 20 comparison-column headers suffixed with
# "_old" / "_new". They are renamed to reference the two concrete file
# format versions being diffed ("_FV2304" / "_FV2310"), the data range is
# turned into a native Excel Table ("Table1"), and the header row is frozen
# so it stays visible while scrolling.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map every old header label (shared-string literal) to its replacement.
$oldToNew = @{
  "Segmentname_old"          = "Segmentname_FV2304"
  "Segmentgruppe_old"        = "Segmentgruppe_FV2304"
  "Segment_old"              = "Segment_FV2304"
  "Datenelement_old"         = "Datenelement_FV2304"
  "Segment ID_old"           = "Segment ID_FV2304"
  "Code_old"                 = "Code_FV2304"
  "Qualifier_old"            = "Qualifier_FV2304"
  "Beschreibung_old"         = "Beschreibung_FV2304"
  "Bedingungsausdruck_old"   = "Bedingungsausdruck_FV2304"
  "Bedingung_old"            = "Bedingung_FV2304"
  "Segmentname_new"          = "Segmentname_FV2310"
  "Segmentgruppe_new"        = "Segmentgruppe_FV2310"
  "Segment_new"              = "Segment_FV2310"
  "Datenelement_new"         = "Datenelement_FV2310"
  "Segment ID_new"           = "Segment ID_FV2310"
  "Code_new"                 = "Code_FV2310"
  "Qualifier_new"            = "Qualifier_FV2310"
  "Beschreibung_new"         = "Beschreibung_FV2310"
  "Bedingungsausdruck_new"   = "Bedingungsausdruck_FV2310"
  "Bedingung_new"            = "Bedingung_FV2310"
}

# The header row is row 1, columns A..U (1..21).
$lastCol = 21
for ($col = 1; $col -le $lastCol; $col++) {
  $cell = $ws.Cells.Item(1, $col)
  $current = $cell.Value()
  if ($oldToNew.ContainsKey($current)) {
    $cell.Value = $oldToNew[$current]
  }
}

# Find the used range so the table covers exactly the data (A1:U75 here).
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
$headerRange = $ws.Range($ws.Cells.Item(1, 1), $ws.Cells.Item($lastRow, $lastCol))

# Turn the range into a proper Excel Table (ListObject) with headers.
$tbl = $ws.ListObjects.Add(1, $headerRange, $null, 1)
$tbl.Name = "Table1"

# Freeze the header row (row 1) so it stays visible while scrolling.
$ws.Activate() | Out-Null
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true
